# The workbook's single data sheet contains one record per row (A:R).
# This commit adds one new weekly record for "Feria Lagunitas de Puerto
# Montt - Ciboulette". The new record is inserted as row 106, which
# pushes the previously-existing rows 106-172 down to rows 107-173
# (dimension grows from A1:R172 to A1:R173).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 106, shifting rows 106:172 down to 107:173.
$ws.Rows("106:106").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A106").Value = 4
$ws.Range("B106").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C106").Value = "Los Lagos"
$ws.Range("D106").Value = 44596
$ws.Range("E106").Value = 10
$ws.Range("F106").Value = 100112039
$ws.Range("G106").Value = "Ciboulette"
$ws.Range("H106").Value = "Sin especificar"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 240
$ws.Range("K106").Value = 2500
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = 2750
$ws.Range("N106").Value = "`$/atado"
$ws.Range("O106").Value = "Provincia de Caut$([char]0x00ED)n"
$ws.Range("P106").Value = 2750
$ws.Range("Q106").Value = 1
$ws.Range("R106").Value = "Hortaliza"

# Keep the same date display/number format as the other "Fecha" cells in
# column D (style index 2 in the original file uses numFmtId 165).
$ws.Range("D106").NumberFormat = $ws.Range("D107").NumberFormat
